$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

$ws.Range("A$row").Value = "8N43SX"
$ws.Range("B$row").Value = "Engranje de transmición del fusor HP"
$ws.Range("C$row").Value = "P1505 P1506 P1566 P1606 M1120 M1522 M1536"
$ws.Range("D$row").Value = 20000
$ws.Range("E$row").Value = 100000
$ws.Range("F$row").Value = 4
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Formula = "=(E$row-D$row)*G$row"
$ws.Range("I$row").Formula = "=D$row*F$row"
$ws.Range("J$row").Value = 80000
